$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cohorts")

# Copy the header formatting from an existing header cell onto the new header cell
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New header for the additional names column
$ws.Range("C1").Value = "Previous/other/additional names"

# New data value for the ABC cohort row
$ws.Range("C2").Value = "ABC-1"
